$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: 199 | LC/GFG | Binary Tree Right Side View
$ws.Range("A11").HorizontalAlignment = -4131
$ws.Range("A11").VerticalAlignment = -4160
$ws.Range("A11").WrapText = $false
$ws.Range("A11").Value = 199
$ws.Range("B11").Value = "LC/GFG"
$ws.Range("C11").Value = "Binary Tree Right Side View"

# Row 12: GFG | GFG | Left View of Binary Tree
$ws.Range("A12").Value = "GFG"
$ws.Range("B12").Value = "GFG"
$ws.Range("C12").Value = "Left View of Binary Tree"

$ws.Range("B11").Select()
